# Helper: populate a new "electricity source" data row (name in col A,
# quantity 2/1 in col B, and the $B<row> fill-across formula in C:AK).
function Add-SourceRow($ws, $r, $name) {
    $ws.Range("A$r").Value = $name
    $ws.Range("B$r").Value = 2
    $ws.Range("C$r`:AK$r").Formula = '=$B' + $r
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDPbES")

# --- Add new electricity-source row 13 (lignite) ---
Add-SourceRow $ws 13 "lignite"

# --- Rename existing electricity-source rows to match new taxonomy ---
$ws.Range("A2").Value = "hard coal"
$ws.Range("A6").Value = "onshore wind"

# --- Append remaining new electricity-source rows (14-17) ---
Add-SourceRow $ws 14 "offshore wind"
Add-SourceRow $ws 15 "crude oil"
Add-SourceRow $ws 16 "heavy or residual fuel oil"
Add-SourceRow $ws 17 "municipal solid waste"

# --- Add new header cell A1 (Priority Order column label) ---
$ws.Range("A1").Value = "Priority Order (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30
